# Fix heat rate modeling syntax
# Applies corrected hourly dispatch / cost values recomputed after fixing
# a heat-rate formula/syntax bug in the underlying model.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("DG Dispatch")
$ws.Range("T2").Value = 0

$ws = $wb.Worksheets.Item("Costs and Revenues")
$ws.Range("B2").Value = 76326.84612799998
$ws.Range("D2").Value = 9307.780929750721
$ws.Range("F2").Value = 41271.41399854876

$ws = $wb.Worksheets.Item("PV Dispatch")
$ws.Range("I3").Value = 41.6
$ws.Range("J3").Value = 11.05311702887466
$ws.Range("M3").Value = 104
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 52
$ws.Range("R3").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 23.4
$ws.Range("N4").Value = 58.18312417100302
$ws.Range("P4").Value = 41.6
$ws.Range("R4").Value = 0

$ws = $wb.Worksheets.Item("Battery Input")
$ws.Range("I3").Value = 41.6
$ws.Range("J3").Value = 11.05311702887466
$ws.Range("M3").Value = 80.59999999999999
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 23.4
$ws.Range("R3").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 58.18312417100301
$ws.Range("P4").Value = 41.6
$ws.Range("R4").Value = 0

$ws = $wb.Worksheets.Item("Battery Output")
$ws.Range("T2").Value = 29.84851999999987

$ws = $wb.Worksheets.Item("State of Charge")
$ws.Range("I3").Value = 171.584
$ws.Range("J3").Value = 182.5265858585859
$ws.Range("K3").Value = 264.8945858585859
$ws.Range("L3").Value = 357.5585858585859
$ws.Range("M3").Value = 437.3525858585859
$ws.Range("N3").Value = 437.3525858585859
$ws.Range("O3").Value = 437.3525858585859
$ws.Range("P3").Value = 460.5185858585859
$ws.Range("Q3").Value = 486.2585858585859
$ws.Range("R3").Value = 486.2585858585859
$ws.Range("L4").Value = 181.88
$ws.Range("M4").Value = 181.88
$ws.Range("N4").Value = 239.481292929293
$ws.Range("O4").Value = 311.553292929293
$ws.Range("P4").Value = 352.7372929292929
$ws.Range("Q4").Value = 373.3292929292929

$ws = $wb.Worksheets.Item("Feed in from Type 2")
$ws.Range("T2").Value = 1.351480000000127

Write-Host "Applied heat rate fix values across DG Dispatch, PV Dispatch, Battery Input, Battery Output, State of Charge, Feed in from Type 2, and Costs and Revenues sheets."
